# Update - Task Schedular.
# Refresh the latest test-run row (row 2) on Sheet1 with the newest
# values recorded by the scheduled test task.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "TestForm_13/03/2019-17:18:25"
$ws.Range("B2").Value = "TestTask_20/02/2019-15:50:19"
$ws.Range("C2").Value = "TestDocument_21/02/2019-15:14:43"
$ws.Range("D2").Value = "TestRSTDocument_21/02/2019-15:14:43"
$ws.Range("E2").Value = "TestAssignTaskDocument21/02/2019-15:14:43"
$ws.Range("F2").Value = "TestDMSTask_21/02/2019-15:14:43"
$ws.Range("I2").Value = "TestForm_05/03/2019-14:55:52"
$ws.Range("J2").Value = "Data Extract-Test2-Form Compliance-2019-02-18-08-57-10.210.xlsx"
